$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 834, pushing the existing data (old rows 834-932)
# down to rows 837-935. This mirrors a new week of price data being prepended
# to the dataset (commit: "Fruta / hortaliza, semanal").
$ws.Rows.Item(834).Resize(3).Insert()

# Populate the 3 newly inserted rows (834-836) with a fresh week of data for
# "Terminal La Palmera de La Serena" / Platano, following the same layout as
# the surrounding rows (now shifted to 837-839).
$ws.Range("A834").Value = 8
$ws.Range("B834").Value = "Terminal La Palmera de La Serena"
$ws.Range("C834").Value = "Coquimbo"
$ws.Range("D834").Value = 44918
$ws.Range("E834").Value = 4
$ws.Range("F834").Value = "Fruta"
$ws.Range("G834").Value = 100108
$ws.Range("H834").Value = "Tropicales y subtropicales"
$ws.Range("I834").Value = 100108006
$ws.Range("J834").Value = "Plátano"
$ws.Range("K834").Value = "Sin especificar"
$ws.Range("L834").Value = "Pintón"
$ws.Range("M834").Value = 80
$ws.Range("N834").Value = 21000
$ws.Range("O834").Value = 21000
$ws.Range("P834").Value = 21000
$ws.Range("Q834").Value = "$/caja 20 kilos"
$ws.Range("R834").Value = "Ecuador"
$ws.Range("S834").Value = 1050
$ws.Range("T834").Value = 20

$ws.Range("A835").Value = 8
$ws.Range("B835").Value = "Terminal La Palmera de La Serena"
$ws.Range("C835").Value = "Coquimbo"
$ws.Range("D835").Value = 44918
$ws.Range("E835").Value = 4
$ws.Range("F835").Value = "Fruta"
$ws.Range("G835").Value = 100108
$ws.Range("H835").Value = "Tropicales y subtropicales"
$ws.Range("I835").Value = 100108006
$ws.Range("J835").Value = "Plátano"
$ws.Range("K835").Value = "Sin especificar"
$ws.Range("L835").Value = "Primera Maduro"
$ws.Range("M835").Value = 120
$ws.Range("N835").Value = 25000
$ws.Range("O835").Value = 25000
$ws.Range("P835").Value = 25000
$ws.Range("Q835").Value = "$/caja 20 kilos"
$ws.Range("R835").Value = "Ecuador"
$ws.Range("S835").Value = 1250
$ws.Range("T835").Value = 20

$ws.Range("A836").Value = 8
$ws.Range("B836").Value = "Terminal La Palmera de La Serena"
$ws.Range("C836").Value = "Coquimbo"
$ws.Range("D836").Value = 44918
$ws.Range("E836").Value = 4
$ws.Range("F836").Value = "Fruta"
$ws.Range("G836").Value = 100108
$ws.Range("H836").Value = "Tropicales y subtropicales"
$ws.Range("I836").Value = 100108006
$ws.Range("J836").Value = "Plátano"
$ws.Range("K836").Value = "Sin especificar"
$ws.Range("L836").Value = "Primera Pintón"
$ws.Range("M836").Value = 120
$ws.Range("N836").Value = 25000
$ws.Range("O836").Value = 25000
$ws.Range("P836").Value = 25000
$ws.Range("Q836").Value = "$/caja 20 kilos"
$ws.Range("R836").Value = "Ecuador"
$ws.Range("S836").Value = 1250
$ws.Range("T836").Value = 20
